$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "60ру" + _GoBack bookmark + "б"  ->  single run "60руб" (bookmark gone)
#    Find/Replace across the whole phrase collapses the three runs (and the
#    bookmark sitting between them) into one run with the merged text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("60руб", $true, $false, $false, $false, $false, $true, 1, $false, "60руб", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Append a page break to the last (empty) paragraph of the body, then add
#    the new "Теплоотвод" section after it, ending with a paragraph that
#    carries the (re-created) _GoBack bookmark.
# ---------------------------------------------------------------------------

# -- 2a: add the page-break run to the existing last paragraph in place.
$endPos = $d.Content.End
$tail = $d.Range($endPos, $endPos)
$breakXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
$tail.InsertXML($breakXml)

# -- 2b: insert the new heading + body paragraphs right after it.
$endPos2 = $d.Content.End
$tail2 = $d.Range($endPos2, $endPos2)
$sectionXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:lastRenderedPageBreak/><w:t>Теплоотвод</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Светодиод:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Тепловое сопротивление </w:t></w:r><w:r><w:t>J</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>-</w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> 10</w:t></w:r><w:r><w:t>K</w:t></w:r><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>/</w:t></w:r><w:r><w:t>W</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="ru-RU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>Максимальная температура активной области +125</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$tail2.InsertXML($sectionXml)
